$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 8699.200000000001
$ws.Range("I18").Value = 8624.25
$ws.Range("K18").Value = 8624.25
$ws.Range("M18").Value = -8340.25

$ws.Range("H69").Value = 18110.846
$ws.Range("I69").Value = 5500
$ws.Range("J69").Value = 20403.727
$ws.Range("K69").Value = 16500
$ws.Range("L69").Value = 61211.181
$ws.Range("M69").Value = -15626
$ws.Range("N69").Value = -62959.181

$ws.Range("H72").Value = 18110.846
$ws.Range("I72").Value = 5500
$ws.Range("J72").Value = 20403.727
$ws.Range("K72").Value = 49500
$ws.Range("L72").Value = 183633.543
$ws.Range("M72").Value = -45132
$ws.Range("N72").Value = -192369.543

$ws.Range("H116").Value = 5000
$ws.Range("I116").Value = 5000
$ws.Range("K116").Value = 5000
$ws.Range("M116").Value = -1558

$ws.Range("H138").Value = 1378
$ws.Range("I138").Value = 1378
$ws.Range("K138").Value = 4134
$ws.Range("M138").Value = 1006

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3180.7222
$ws.Range("I32").Value = 3180.7222
$ws.Range("K32").Value = 3180.7222
$ws.Range("M32").Value = -2893.7222

$ws.Range("H45").Value = 6499.5
$ws.Range("I45").Value = 6499.5
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 6499.5
$ws.Range("L45").Value = 0
$ws.Range("N45").Value = -6122.5
$ws.Range("M45").ClearContents()

$ws.Range("H61").Value = 3140.0527
$ws.Range("I61").Value = 3003.9375
$ws.Range("K61").Value = 3003.9375
$ws.Range("M61").Value = -2791.9375

$ws.Range("H74").Value = 1422.0555
$ws.Range("I74").Value = 1281.6875
$ws.Range("K74").Value = 1281.6875
$ws.Range("M74").Value = -407.6875

$ws.Range("H77").Value = 1422.0555
$ws.Range("I77").Value = 1281.6875
$ws.Range("K77").Value = 6408.4375
$ws.Range("M77").Value = -2040.4375

$ws.Range("H102").Value = 363
$ws.Range("I102").Value = 363
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 363
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = 1259
$ws.Range("M102").ClearContents()

$ws.Range("H132").Value = 2256.3125
$ws.Range("I132").Value = 2256.3125
$ws.Range("K132").Value = 6768.9375
$ws.Range("M132").Value = -4238.9375

$ws.Range("H136").Value = 3140.0527
$ws.Range("I136").Value = 3003.9375
$ws.Range("K136").Value = 9011.8125
$ws.Range("M136").Value = -6461.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 45000
$ws.Range("J92").Value = 45000
$ws.Range("L92").Value = 45000
$ws.Range("N92").Value = -49992

$ws.Range("H99").Value = 1638.5555
$ws.Range("I99").Value = 1193.6
$ws.Range("J99").Value = 2194.75
$ws.Range("K99").Value = 1193.6
$ws.Range("L99").Value = 2194.75
$ws.Range("M99").Value = 304.4000000000001
$ws.Range("N99").Value = -5190.75

$ws.Range("H134").Value = 5457.4414
$ws.Range("I134").Value = 5219.0347
$ws.Range("J134").Value = 6840.2
$ws.Range("K134").Value = 15657.1041
$ws.Range("L134").Value = 20520.6
$ws.Range("M134").Value = -13122.1041
$ws.Range("N134").Value = -25590.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2250
$ws.Range("J16").Value = 3000
$ws.Range("L16").Value = 3000
$ws.Range("N16").Value = -3574

$ws.Range("H22").Value = 8000519
$ws.Range("J22").Value = 10000600
$ws.Range("L22").Value = 10000600
$ws.Range("N22").Value = -10001300

$ws.Range("H86").Value = 9993.4
$ws.Range("I86").Value = 9991.75
$ws.Range("K86").Value = 9991.75
$ws.Range("M86").Value = -8868.75

$ws.Range("H89").Value = 9993.4
$ws.Range("I89").Value = 9991.75
$ws.Range("K89").Value = 49958.75
$ws.Range("M89").Value = -44342.75

$ws.Range("H94").Value = 1494.6666
$ws.Range("I94").Value = 1393.8
$ws.Range("J94").Value = 1999
$ws.Range("K94").Value = 1393.8
$ws.Range("L94").Value = 1999
$ws.Range("M94").Value = -942.8
$ws.Range("N94").Value = -2901

$ws.Range("H113").Value = 2250
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340

$ws.Range("H132").Value = 1639.6
$ws.Range("I132").Value = 1624.5
$ws.Range("J132").Value = 1700
$ws.Range("K132").Value = 4873.5
$ws.Range("L132").Value = 5100
$ws.Range("M132").Value = -2343.5
$ws.Range("N132").Value = -10160

$ws.Range("H134").Value = 1844.5588
$ws.Range("I134").Value = 1839.8485
$ws.Range("K134").Value = 5519.5455
$ws.Range("M134").Value = -2984.5455

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 111.14286
$ws.Range("I6").Value = 88
$ws.Range("K6").Value = 264
$ws.Range("M6").Value = -151

$ws.Range("H11").Value = 218.2
$ws.Range("I11").Value = 218.2
$ws.Range("K11").Value = 654.5999999999999
$ws.Range("M11").Value = -514.5999999999999

$ws.Range("H38").Value = 279.75
$ws.Range("I38").Value = 279.75
$ws.Range("K38").Value = 839.25
$ws.Range("M38").Value = -492.25

$ws.Range("I107").Value = 999
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 2997
$ws.Range("L107").Value = 2700
$ws.Range("M107").Value = -1077
$ws.Range("N107").Value = -6540

$ws.Range("H109").Value = 4960.357
$ws.Range("I109").Value = 1064
$ws.Range("K109").Value = 3192
$ws.Range("M109").Value = -2152

$ws.Range("H129").Value = 2708.25
$ws.Range("I129").Value = 800
$ws.Range("J129").Value = 3344.3333
$ws.Range("K129").Value = 2400
$ws.Range("L129").Value = 10032.9999
$ws.Range("M129").Value = 2600
$ws.Range("N129").Value = -20032.9999

$ws.Range("H131").Value = 1946.758
$ws.Range("I131").Value = 1900
$ws.Range("K131").Value = 5700
$ws.Range("M131").Value = -660

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 9500
$ws.Range("K70").Value = 9500
$ws.Range("M70").Value = -9230

$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 9500
$ws.Range("K73").Value = 9500
$ws.Range("M73").Value = -8564

$ws.Range("H102").Value = 1485.6
$ws.Range("I102").Value = 1485.6
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1485.6
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = 136.4000000000001
$ws.Range("M102").ClearContents()

$ws.Range("H122").Value = 4286.5386
$ws.Range("I122").Value = 4516.5
$ws.Range("J122").Value = 3918.6
$ws.Range("K122").Value = 13549.5
$ws.Range("L122").Value = 11755.8
$ws.Range("M122").Value = -11099.5
$ws.Range("N122").Value = -16655.8

$ws.Range("H132").Value = 2396.7
$ws.Range("I132").Value = 2223.8572
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 6671.571599999999
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = -4141.571599999999
$ws.Range("N132").Value = -13460

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 473.75
$ws.Range("J9").Value = 149
$ws.Range("L9").Value = 149
$ws.Range("N9").Value = -597

$ws.Range("H122").Value = 2581.3333
$ws.Range("I122").Value = 2581.3333
$ws.Range("K122").Value = 7743.999899999999
$ws.Range("M122").Value = -5293.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 67499
$ws.Range("J125").Value = 67499
$ws.Range("L125").Value = 67499
$ws.Range("N125").Value = -77339

$ws.Range("H126").Value = 5194.1665
$ws.Range("I126").Value = 4370.4443
$ws.Range("K126").Value = 13111.3329
$ws.Range("M126").Value = -10641.3329
